$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 6.353900000000002
$ws.Range("B21").Value = 9.567400000000006
$ws.Range("B23").Value = 8.909500000000005
$ws.Range("B25").Value = 5.857300000000001
$ws.Range("B53").Value = 5.6369
$ws.Range("B57").Value = 4.993699999999997
$ws.Range("B59").Value = 4.812099999999996
$ws.Range("B69").Value = 5.394199999999996
$ws.Range("B79").Value = 8.650300000000001
$ws.Range("B83").Value = 5.403999999999999
$ws.Range("B93").Value = 5.7704
